# updated PvsI with model fitting
# Applies refreshed volume/rate calculations (columns T, Z, AB, AD) for rows 10-17
# on the active worksheet, matching a re-run of the respirometry rate-fitting model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 10; T = 0.1405756097560976; Z = 0.5226268671296774;  AB = 2108.663240577803;  AD = 2108.663240577803 },
    @{ Row = 11; T = 0.144009756097561;  Z = 0.1435213970517591;  AB = 945.4843242616821;  AD = 945.4843242616821 },
    @{ Row = 12; T = 0.1462634146341464; Z = 0.05872983817896857; AB = 376.2279544972676;  AD = 376.2279544972676 },
    @{ Row = 13; T = 0.1500390243902439; Z = 0.09971952446600377; AB = 622.775104237067;   AD = 622.775104237067 },
    @{ Row = 14; T = 0.1463707317073171; Z = 0.03825761303729501; AB = 238.9945947258462;  AD = 238.9945947258462 },
    @{ Row = 15; T = 0.1465658536585366; Z = 0.1269890575083733;  AB = 633.6648191847016;  AD = 633.6648191847016 },
    @{ Row = 16; T = 0.1449658536585366; Z = 0.1849666055079211;  AB = 746.093702096627;   AD = 746.093702096627 },
    @{ Row = 17; T = 0.1544 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("T$r").Value = $u.T
    if ($u.ContainsKey("Z"))  { $ws.Range("Z$r").Value  = $u.Z }
    if ($u.ContainsKey("AB")) { $ws.Range("AB$r").Value = $u.AB }
    if ($u.ContainsKey("AD")) { $ws.Range("AD$r").Value = $u.AD }
}
